# Added logic for Runmode in testdata xls in listener onTestStart()
# Adds three new rows of test data (Rahul, Ishita, Rohit Sehgal) to the
# "AddCustomerTest" sheet, mirroring the existing Raman/Arora row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddCustomerTest")

$ws.Range("A3").Value = "Rahul"
$ws.Range("B3").Value = "Arora"
$ws.Range("C3").Value = "A234wd"
$ws.Range("D3").Value = "Customer added successfully"

$ws.Range("A4").Value = "Ishita"
$ws.Range("B4").Value = "Arora"
$ws.Range("C4").Value = "A234wd"
$ws.Range("D4").Value = "Customer added successfully"

$ws.Range("A5").Value = "Rohit"
$ws.Range("B5").Value = "Sehgal"
$ws.Range("C5").Value = "A234wd"
$ws.Range("D5").Value = "Customer added successfully"

$ws.Activate()
$ws.Range("B5").Select()
